# Trade #100 closed at 2026-02-17 21:27:33 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Summary
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1401.39
$summary.Range("B4").Value = 1.18
$summary.Range("B5").Value = 0.18
$summary.Range("B6").Value = 128
$summary.Range("B8").Value = 49
$summary.Range("B9").Value = 42.97

# ---------------------------------------------------------------------------
# Sheet: Strategy Status (row 5 = MarketMaking)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 101.39
$status.Range("D5").Value = 95
$status.Range("E5").Value = 1.07
$status.Range("F5").Value = 1.39
$status.Range("G5").Value = 43.16

# ---------------------------------------------------------------------------
# Sheet: All Trades (Trade #128 is in row 129; close it out)
# Columns: A Trade#, B Date, C Time, D Strategy, E Side, F Entry Price,
#          G Exit Price, H Status, I P&L %, J P&L $, K Capital After,
#          L Exit Reason, M Duration (min), N Entry Slippage (bps),
#          O Exit Slippage (bps), P Confidence, Q Entry Reason
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Cells.Item(129, 7).Value = 0.08823499999999999
$allTrades.Cells.Item(129, 8).Value = "CLOSED"
$allTrades.Cells.Item(129, 9).Value = -11.7647
$allTrades.Cells.Item(129, 10).Value = -0.01
$allTrades.Cells.Item(129, 11).Value = 101.39
$allTrades.Cells.Item(129, 12).Value = "early_exit"
$allTrades.Cells.Item(129, 13).Value = 0.13

# New open trade (#161) recorded as row 162
$allTrades.Cells.Item(162, 1).Value = 161
# Copy an existing "2026-02-17" text cell so the date-formatted string is not
# re-interpreted as a date serial number by the COM layer.
$allTrades.Range("B2").Copy($allTrades.Range("B162"))
$allTrades.Cells.Item(162, 3).Value = "21:27:27"
$allTrades.Cells.Item(162, 4).Value = "MarketMaking"
$allTrades.Cells.Item(162, 5).Value = "UP"
$allTrades.Cells.Item(162, 6).Value = 0.1
$allTrades.Cells.Item(162, 8).Value = "OPEN"
$allTrades.Cells.Item(162, 9).Value = 0
$allTrades.Cells.Item(162, 10).Value = 0
$allTrades.Cells.Item(162, 11).Value = 101.401797784678
$allTrades.Cells.Item(162, 13).Value = 0
$allTrades.Cells.Item(162, 14).Value = 0
$allTrades.Cells.Item(162, 15).Value = 0
$allTrades.Cells.Item(162, 16).Value = 0.6
$allTrades.Cells.Item(162, 17).Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# Sheet: MarketMaking (Trade #128 is in row 96; close it out)
# Columns: A Trade#, B Date, C Time, D Strategy, E Side, F Entry Price,
#          G Exit Price, H Status, I P&L %, J P&L $, K Capital After,
#          L Entry Slippage (bps), M Exit Slippage (bps), N Confidence,
#          O Entry Reason, P Exit Reason, Q Duration (min)
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Cells.Item(96, 7).Value = 0.08823499999999999
$mm.Cells.Item(96, 8).Value = "CLOSED"
$mm.Cells.Item(96, 9).Value = -11.7647
$mm.Cells.Item(96, 10).Value = -0.01
$mm.Cells.Item(96, 11).Value = 101.39
$mm.Cells.Item(96, 16).Value = "early_exit"
$mm.Cells.Item(96, 17).Value = 0.13

# New open trade (#161) recorded as row 129
$mm.Cells.Item(129, 1).Value = 161
# Copy an existing "2026-02-17" text cell so the date-formatted string is not
# re-interpreted as a date serial number by the COM layer.
$mm.Range("B2").Copy($mm.Range("B129"))
$mm.Cells.Item(129, 3).Value = "21:27:27"
$mm.Cells.Item(129, 4).Value = "MarketMaking"
$mm.Cells.Item(129, 5).Value = "UP"
$mm.Cells.Item(129, 6).Value = 0.1
$mm.Cells.Item(129, 8).Value = "OPEN"
$mm.Cells.Item(129, 9).Value = 0
$mm.Cells.Item(129, 10).Value = 0
$mm.Cells.Item(129, 11).Value = 101.401797784678
$mm.Cells.Item(129, 12).Value = 0
$mm.Cells.Item(129, 13).Value = 0
$mm.Cells.Item(129, 14).Value = 0.6
$mm.Cells.Item(129, 15).Value = "Normal spread capture: 19600 bps"
$mm.Cells.Item(129, 17).Value = 0
